# "se actualiza la data para pagos"
# Update the payment data on the "CronogramaPagosVehicular" sheet and make it
# the active sheet/selection (matching the commit's data refresh).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("CronogramaPagosVehicular")

# Update row 2 values, keeping them as text (the original cells use
# quote-prefixed text styles, even the "date" column), via a leading
# apostrophe so Excel stores them as shared-string text rather than
# auto-converting to number/date.
$ws1.Range("A2").Value = "'72636759"
$ws1.Range("C2").Value = "'ZIJ-583"
$ws1.Range("D2").Value = "'12/11/2025"

# Move the selection to D2 and make this sheet the active one.
$ws1.Range("D2").Select()
$ws1.Activate()
